$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-06-15 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-16 Sunday", 2) | Out-Null

# Update the division-problem answers in the table, addressed by (row, column)
# so that identical/overlapping text values across cells do not collide.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "59÷5=11, 4"
$tbl.Cell(1, 2).Range.Text = "37÷2=18, 1"
$tbl.Cell(1, 3).Range.Text = "82÷5=16, 2"
$tbl.Cell(1, 4).Range.Text = "19÷8=2, 3"
$tbl.Cell(1, 5).Range.Text = "20÷4=5, 0"
$tbl.Cell(5, 1).Range.Text = "64÷2=32, 0"
$tbl.Cell(5, 2).Range.Text = "48÷5=9, 3"
$tbl.Cell(5, 3).Range.Text = "32÷3=10, 2"
$tbl.Cell(5, 4).Range.Text = "92÷9=10, 2"
$tbl.Cell(5, 5).Range.Text = "82÷8=10, 2"
$tbl.Cell(9, 1).Range.Text = "17÷5=3, 2"
$tbl.Cell(9, 2).Range.Text = "51÷3=17, 0"
$tbl.Cell(9, 3).Range.Text = "18÷8=2, 2"
$tbl.Cell(9, 4).Range.Text = "95÷3=31, 2"
$tbl.Cell(9, 5).Range.Text = "34÷8=4, 2"
$tbl.Cell(13, 1).Range.Text = "23÷8=2, 7"
$tbl.Cell(13, 2).Range.Text = "65÷9=7, 2"
$tbl.Cell(13, 3).Range.Text = "46÷2=23, 0"
$tbl.Cell(13, 4).Range.Text = "92÷2=46, 0"
$tbl.Cell(13, 5).Range.Text = "94÷6=15, 4"
$tbl.Cell(17, 1).Range.Text = "13÷4=3, 1"
$tbl.Cell(17, 2).Range.Text = "20÷9=2, 2"
$tbl.Cell(17, 3).Range.Text = "10÷2=5, 0"
$tbl.Cell(17, 4).Range.Text = "29÷9=3, 2"
$tbl.Cell(17, 5).Range.Text = "95÷6=15, 5"
